# Adding the changes we made on may 9th
# Re-selects a shifted 30-sample window of accelerometer/gyroscope readings
# (9 brand-new leading rows + 1 brand-new trailing row, with the previously
# selected rows shifting down by 9) and extends the sheet from 20 to 30
# data rows (A1:H21 -> A1:H31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full replacement block (30 rows x 8 cols: timestamp, label, ax, ay, az, gx, gy, gz)
# as a true 2-D array -- Range.Value needs this shape (not a PS jagged array) to
# write a rectangular block in one shot.
$arr = New-Object 'object[,]' 30,8
$arr[0,0] = 0
$arr[0,1] = "struggle"
$arr[0,2] = -0.188694953918457
$arr[0,3] = -0.0127399563789367
$arr[0,4] = 0.0153613984584808
$arr[0,5] = -0.0114537235349416
$arr[0,6] = -0.0096211275085806
$arr[0,7] = -0.0482583530247211
$arr[1,0] = 100
$arr[1,1] = "struggle"
$arr[1,2] = -0.0261173248291015
$arr[1,3] = -0.1474769711494445
$arr[1,4] = 0.0655251443386077
$arr[1,5] = 0.0445931628346443
$arr[1,6] = 0.1122464910149574
$arr[1,7] = -0.0378736443817615
$arr[2,0] = 200
$arr[2,1] = "struggle"
$arr[2,2] = -0.1960973739624023
$arr[2,3] = 0.0549294650554657
$arr[2,4] = 0.0360765755176544
$arr[2,5] = 0.0612392425537109
$arr[2,6] = 0.09758572280406951
$arr[2,7] = -0.0021380283869802
$arr[3,0] = 300
$arr[3,1] = "struggle"
$arr[3,2] = -0.06610202789306641
$arr[3,3] = -0.1787786185741424
$arr[3,4] = 0.0745508223772049
$arr[3,5] = 0.0088575463742017
$arr[3,6] = 0.1237002089619636
$arr[3,7] = 0.0548251569271087
$arr[4,0] = 400
$arr[4,1] = "struggle"
$arr[4,2] = 0.0234136581420898
$arr[4,3] = 0.0270741879940032
$arr[4,4] = 0.2239813506603241
$arr[4,5] = -0.0221438650041818
$arr[4,6] = 0.0061086523346602
$arr[4,7] = 0.0325285755097866
$arr[5,0] = 500
$arr[5,1] = "struggle"
$arr[5,2] = 0.11651611328125
$arr[5,3] = -0.4856438636779785
$arr[5,4] = 0.5658785104751587
$arr[5,5] = 0.0332921557128429
$arr[5,6] = -0.0615446716547012
$arr[5,7] = 0.093156948685646
$arr[6,0] = 600
$arr[6,1] = "struggle"
$arr[6,2] = 0.0557413101196289
$arr[6,3] = 0.3574482798576355
$arr[6,4] = 0.2321825623512268
$arr[6,5] = -0.4489859640598297
$arr[6,6] = -1.353219270706177
$arr[6,7] = 0.4497495293617248
$arr[7,0] = 700
$arr[7,1] = "struggle"
$arr[7,2] = 0.3619680404663086
$arr[7,3] = 0.0124948024749755
$arr[7,4] = 0.3587799966335296
$arr[7,5] = -0.3888157308101654
$arr[7,6] = -3.63083028793335
$arr[7,7] = -0.1369865238666534
$arr[8,0] = 800
$arr[8,1] = "struggle"
$arr[8,2] = -0.2529764175415039
$arr[8,3] = 0.1160029470920562
$arr[8,4] = -0.09882223606109609
$arr[8,5] = -0.6565274000167847
$arr[8,6] = -2.371837139129639
$arr[8,7] = 0.1600466966629028
$arr[9,0] = 900
$arr[9,1] = "struggle"
$arr[9,2] = -0.1584005355834961
$arr[9,3] = 0.0559865832328796
$arr[9,4] = -0.2031860947608947
$arr[9,5] = -0.4257730841636657
$arr[9,6] = -1.438740372657776
$arr[9,7] = 0.180816113948822
$arr[10,0] = 1000
$arr[10,1] = "struggle"
$arr[10,2] = -0.1681756973266601
$arr[10,3] = -0.045459896326065
$arr[10,4] = 0.3079473972320556
$arr[10,5] = -0.2063197344541549
$arr[10,6] = 0.5047274231910706
$arr[10,7] = -0.1090394482016563
$arr[11,0] = 1100
$arr[11,1] = "struggle"
$arr[11,2] = 0.7375173568725586
$arr[11,3] = -0.8549392819404602
$arr[11,4] = -2.997310400009156
$arr[11,5] = 1.435685992240906
$arr[11,6] = 5.099197864532471
$arr[11,7] = -0.6409503817558289
$arr[12,0] = 1200
$arr[12,1] = "struggle"
$arr[12,2] = -0.6316938400268555
$arr[12,3] = 0.0533061251044273
$arr[12,4] = -1.823783159255981
$arr[12,5] = 0.7269296646118164
$arr[12,6] = 4.458247184753418
$arr[12,7] = 0.2814561724662781
$arr[13,0] = 1300
$arr[13,1] = "struggle"
$arr[13,2] = 0.1245284080505371
$arr[13,3] = 0.4134435057640075
$arr[13,4] = 2.055456638336182
$arr[13,5] = 0.2370157092809677
$arr[13,6] = 0.7996225953102112
$arr[13,7] = 0.1328631937503814
$arr[14,0] = 1400
$arr[14,1] = "struggle"
$arr[14,2] = -1.905292510986328
$arr[14,3] = 1.267569422721863
$arr[14,4] = 0.3008813858032226
$arr[14,5] = 0.2102903574705124
$arr[14,6] = 1.452026724815369
$arr[14,7] = 0.2237294018268585
$arr[15,0] = 1500
$arr[15,1] = "struggle"
$arr[15,2] = -2.353589773178101
$arr[15,3] = 0.5766786336898804
$arr[15,4] = 2.404436111450196
$arr[15,5] = -0.3715587854385376
$arr[15,6] = 0.4751004576683044
$arr[15,7] = 0.1111774742603302
$arr[16,0] = 1600
$arr[16,1] = "struggle"
$arr[16,2] = -11.09067344665527
$arr[16,3] = 1.405970811843872
$arr[16,4] = 10.02403450012207
$arr[16,5] = 0.1693623960018158
$arr[16,6] = -1.752572417259216
$arr[16,7] = 0.1539380401372909
$arr[17,0] = 1700
$arr[17,1] = "struggle"
$arr[17,2] = 4.286171913146973
$arr[17,3] = 0.2758489847183227
$arr[17,4] = -4.509784698486328
$arr[17,5] = -1.307862520217896
$arr[17,6] = -5.349499702453613
$arr[17,7] = -1.575574159622192
$arr[18,0] = 1800
$arr[18,1] = "struggle"
$arr[18,2] = -1.000519752502441
$arr[18,3] = -0.010628342628479
$arr[18,4] = -1.670511245727539
$arr[18,5] = -0.3645338416099548
$arr[18,6] = -2.762179851531982
$arr[18,7] = 0.608421802520752
$arr[19,0] = 1900
$arr[19,1] = "struggle"
$arr[19,2] = -2.810617446899414
$arr[19,3] = 0.8466755151748657
$arr[19,4] = -0.6261429786682129
$arr[19,5] = -0.3593414723873138
$arr[19,6] = -2.416888236999512
$arr[19,7] = -0.4506658315658569
$arr[20,0] = 2000
$arr[20,1] = "struggle"
$arr[20,2] = -0.552617073059082
$arr[20,3] = 1.007189750671387
$arr[20,4] = -2.683732509613037
$arr[20,5] = -0.3178026378154754
$arr[20,6] = -1.223715782165527
$arr[20,7] = -0.2168571650981903
$arr[21,0] = 2100
$arr[21,1] = "struggle"
$arr[21,2] = -2.832679748535156
$arr[21,3] = 5.107204437255859
$arr[21,4] = -6.522222995758057
$arr[21,5] = 0.042302418500185
$arr[21,6] = 0.5458080768585205
$arr[21,7] = 0.195171445608139
$arr[22,0] = 2200
$arr[22,1] = "struggle"
$arr[22,2] = 0.8469958305358887
$arr[22,3] = -1.08077871799469
$arr[22,4] = 7.442714691162109
$arr[22,5] = 1.255175352096558
$arr[22,6] = 4.058435916900635
$arr[22,7] = 0.6265950202941895
$arr[23,0] = 2300
$arr[23,1] = "struggle"
$arr[23,2] = -3.03963303565979
$arr[23,3] = 1.802032470703125
$arr[23,4] = -2.227274417877197
$arr[23,5] = 1.706299304962158
$arr[23,6] = 4.895015716552734
$arr[23,7] = -0.6637051105499268
$arr[24,0] = 2400
$arr[24,1] = "struggle"
$arr[24,2] = -1.961796522140503
$arr[24,3] = 1.68219518661499
$arr[24,4] = 1.394426345825195
$arr[24,5] = 0.3026837408542633
$arr[24,6] = 0.5484042763710022
$arr[24,7] = 0.1058324053883552
$arr[25,0] = 2500
$arr[25,1] = "struggle"
$arr[25,2] = -2.372189998626709
$arr[25,3] = 1.225671410560608
$arr[25,4] = 2.504203796386719
$arr[25,5] = -0.2704605758190155
$arr[25,6] = 0.6455318331718445
$arr[25,7] = -0.1504255682229995
$arr[26,0] = 2600
$arr[26,1] = "struggle"
$arr[26,2] = -8.473310470581055
$arr[26,3] = -0.7327957153320312
$arr[26,4] = 5.200639724731445
$arr[26,5] = -0.4751004576683044
$arr[26,6] = -0.1533271819353103
$arr[26,7] = -0.3729332387447357
$arr[27,0] = 2700
$arr[27,1] = "struggle"
$arr[27,2] = 3.173869132995605
$arr[27,3] = -1.535699486732483
$arr[27,4] = -6.114311695098877
$arr[27,5] = -0.3394883573055267
$arr[27,6] = -0.6884451508522034
$arr[27,7] = -0.1850921660661697
$arr[28,0] = 2800
$arr[28,1] = "struggle"
$arr[28,2] = -4.002721786499023
$arr[28,3] = 1.022015571594239
$arr[28,4] = -0.0432633161544799
$arr[28,5] = -0.2115120887756347
$arr[28,6] = -0.2267837226390838
$arr[28,7] = -0.3090978264808655
$arr[29,0] = 2900
$arr[29,1] = "struggle"
$arr[29,2] = -0.8564167022705078
$arr[29,3] = -0.1756476759910583
$arr[29,4] = -1.401212096214294
$arr[29,5] = 0.0546724386513233
$arr[29,6] = -0.1007927656173706
$arr[29,7] = 0.2141082733869552

$ws.Range("A2:H31").Value = $arr

Write-Output "Updated Sheet1!A2:H31 (30 data rows)"
